$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "missing data" simulation was re-rolled: two rows (RM 232 and SC 92)
# that used to be included are no longer part of the sample, and a
# different set of cells are now blanked out / filled in to represent
# which values are treated as missing.

# Remove the "RM 232" row (old row 26). Rows below shift up by one.
$ws.Rows("26").Delete()
# Remove the "SC 92" row (old row 28, now row 27 after the first delete).
$ws.Rows("27").Delete()

# Helper: blank out a cell the same way the source data does it -- an
# empty-text cell (not a fully-blank/non-existent cell). Using an
# apostrophe forces Excel to store a literal empty string, then we reset
# the style so no stray "quote prefix" formatting is left behind.
function Set-Missing($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

# F5 (RM 14) used to have a value, now it's missing.
Set-Missing "F5"

# F8 (RM 38) used to be missing, now it has a value.
$ws.Range("F8").Value = 17.05

# F12 (RM 81) used to have a value, now it's missing.
Set-Missing "F12"

# F14 (RM 90) used to be missing, now it has a value.
$ws.Range("F14").Value = 17.76

# F18 (RM 120) used to have a value, now it's missing.
Set-Missing "F18"

# B26 (SC 5, after row shift) used to be missing, now it has a value.
$ws.Range("B26").Value = -20.2

# B27 (SC 101, after row shift) used to have a value, now it's missing.
Set-Missing "B27"

# E33 (SC 232, after row shift) used to be missing, now it has a value.
$ws.Range("E33").Value = -10.7
